$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 212.14285
$ws.Range("I2").Value = 97
$ws.Range("K2").Value = 97
$ws.Range("M2").Value = 16
$ws.Range("H6").Value = 77.71429000000001
$ws.Range("I6").Value = 57.333332
$ws.Range("K6").Value = 171.999996
$ws.Range("M6").Value = -59.99999600000001
$ws.Range("H9").Value = 660.7273
$ws.Range("I9").Value = 420.875
$ws.Range("K9").Value = 420.875
$ws.Range("M9").Value = -251.875
$ws.Range("H40").Value = 2400
$ws.Range("J40").Value = 2400
$ws.Range("L40").Value = 2400
$ws.Range("N40").Value = -2750
$ws.Range("H86").Value = 24256.223
$ws.Range("I86").Value = 20320.2
$ws.Range("K86").Value = 20320.2
$ws.Range("M86").Value = -19197.2
$ws.Range("H89").Value = 24256.223
$ws.Range("I89").Value = 20320.2
$ws.Range("K89").Value = 101601
$ws.Range("M89").Value = -95985
$ws.Range("H96").Value = 2176.4614
$ws.Range("I96").Value = 2154.5
$ws.Range("J96").Value = 2249.6667
$ws.Range("K96").Value = 6463.5
$ws.Range("L96").Value = 6749.000100000001
$ws.Range("M96").Value = -5090.5
$ws.Range("N96").Value = -9495.000100000001
$ws.Range("H98").Value = 5828.1113
$ws.Range("I98").Value = 5828.1113
$ws.Range("K98").Value = 5828.1113
$ws.Range("M98").Value = -4330.1113
$ws.Range("H111").Value = 531.6842
$ws.Range("I111").Value = 301
$ws.Range("K111").Value = 903
$ws.Range("M111").Value = 2164
$ws.Range("H122").Value = 5828.1113
$ws.Range("I122").Value = 5828.1113
$ws.Range("K122").Value = 17484.3339
$ws.Range("M122").Value = -15034.3339
$ws.Range("H127").Value = 357967.28
$ws.Range("I127").Value = 357967.28
$ws.Range("K127").Value = 1073901.84
$ws.Range("M127").Value = -1068941.84
$ws.Range("H132").Value = 10998.429
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
$ws.Range("H138").Value = 3650.7795
$ws.Range("I138").Value = 3307.476
$ws.Range("K138").Value = 9922.428
$ws.Range("M138").Value = -4782.428

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1851123.2
$ws.Range("I61").Value = 8049.5
$ws.Range("K61").Value = 8049.5
$ws.Range("M61").Value = -7837.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 5765109.5
$ws.Range("I132").Value = 7637.2354
$ws.Range("K132").Value = 22911.7062
$ws.Range("M132").Value = -20381.7062
$ws.Range("H136").Value = 1851123.2
$ws.Range("I136").Value = 8049.5
$ws.Range("K136").Value = 24148.5
$ws.Range("M136").Value = -21598.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H96").Value = 24166.334
$ws.Range("I96").Value = 24166.334
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 24166.334
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -21420.334
$ws.Range("N96").ClearContents()
$ws.Range("H105").Value = 1946.9286
$ws.Range("J105").Value = 1905
$ws.Range("L105").Value = 1905
$ws.Range("N105").Value = -5399

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19009.834
$ws.Range("I58").Value = 9323.375
$ws.Range("K58").Value = 9323.375
$ws.Range("M58").Value = -9120.375
$ws.Range("H132").Value = 33213204
$ws.Range("I132").Value = 3481.5
$ws.Range("K132").Value = 10444.5
$ws.Range("M132").Value = -7914.5
$ws.Range("H136").Value = 19009.834
$ws.Range("I136").Value = 9323.375
$ws.Range("K136").Value = 27970.125
$ws.Range("M136").Value = -25420.125
$ws.Range("H140").Value = 171187
$ws.Range("J140").Value = 165839.67
$ws.Range("L140").Value = 165839.67
$ws.Range("N140").Value = -176199.67

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3724.4614
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 3724.4614
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 11173.3842
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -11453.3842
$ws.Range("H38").Value = 86.08
$ws.Range("I38").Value = 89.21429000000001
$ws.Range("J38").Value = 82.09090999999999
$ws.Range("K38").Value = 267.64287
$ws.Range("L38").Value = 246.27273
$ws.Range("M38").Value = 79.35712999999998
$ws.Range("N38").Value = -940.2727299999999
$ws.Range("H68").Value = 1481.55
$ws.Range("J68").Value = 1481.55
$ws.Range("L68").Value = 4444.65
$ws.Range("N68").Value = -6066.65
$ws.Range("H71").Value = 1481.55
$ws.Range("J71").Value = 1481.55
$ws.Range("L71").Value = 13333.95
$ws.Range("N71").Value = -21445.95
$ws.Range("H121").Value = 6471987.5
$ws.Range("I121").Value = 1667117.4
$ws.Range("J121").Value = 9092825
$ws.Range("K121").Value = 5001352.199999999
$ws.Range("L121").Value = 27278475
$ws.Range("M121").Value = -5000042.199999999
$ws.Range("N121").Value = -27281095
$ws.Range("H131").Value = 1456.5918
$ws.Range("J131").Value = 1473.5745
$ws.Range("L131").Value = 4420.7235
$ws.Range("N131").Value = -14500.7235
$ws.Range("H132").Value = 1685106.6
$ws.Range("I132").Value = 1554
$ws.Range("K132").Value = 13986
$ws.Range("M132").Value = -11456

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5619
$ws.Range("I102").Value = 6253.3076
$ws.Range("K102").Value = 6253.3076
$ws.Range("M102").Value = -4631.3076
$ws.Range("H126").Value = 9973.875
$ws.Range("I126").Value = 8272.166999999999
$ws.Range("K126").Value = 24816.501
$ws.Range("M126").Value = -22346.501

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3274
$ws.Range("I16").Value = 3274
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3274
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3104
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 4399.5
$ws.Range("J22").Value = 4399.5
$ws.Range("L22").Value = 4399.5
$ws.Range("N22").Value = -4989.5
$ws.Range("H27").Value = 4399.5
$ws.Range("J27").Value = 4399.5
$ws.Range("L27").Value = 4399.5
$ws.Range("N27").Value = -4613.5
$ws.Range("H40").Value = 4152.278
$ws.Range("I40").Value = 2874.5417
$ws.Range("J40").Value = 6707.75
$ws.Range("K40").Value = 2874.5417
$ws.Range("L40").Value = 6707.75
$ws.Range("M40").Value = -2738.5417
$ws.Range("N40").Value = -6979.75
$ws.Range("H55").Value = 573.03705
$ws.Range("J55").Value = 722.36365
$ws.Range("L55").Value = 722.36365
$ws.Range("N55").Value = -1068.36365
$ws.Range("H61").Value = 3429.4849
$ws.Range("I61").Value = 2623.077
$ws.Range("K61").Value = 2623.077
$ws.Range("M61").Value = -2421.077
$ws.Range("H96").Value = 16666.666
$ws.Range("J96").Value = 16666.666
$ws.Range("L96").Value = 16666.666
$ws.Range("N96").Value = -22158.666
$ws.Range("H113").Value = 3429.4849
$ws.Range("I113").Value = 2623.077
$ws.Range("K113").Value = 2623.077
$ws.Range("M113").Value = -453.0770000000002
$ws.Range("H122").Value = 7294.5
$ws.Range("I122").Value = 7464
$ws.Range("K122").Value = 22392
$ws.Range("M122").Value = -19942
$ws.Range("H132").Value = 1593705.8
$ws.Range("I132").Value = 4326.8096
$ws.Range("K132").Value = 12980.4288
$ws.Range("M132").Value = -10450.4288
$ws.Range("H136").Value = 700119.4
$ws.Range("I136").Value = 18915
$ws.Range("K136").Value = 56745
$ws.Range("M136").Value = -54195

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2553686
$ws.Range("J132").Value = 4886851.5
$ws.Range("L132").Value = 14660554.5
$ws.Range("N132").Value = -14665614.5
$ws.Range("H136").Value = 821420.9
$ws.Range("J136").Value = 889122.5600000001
$ws.Range("L136").Value = 2667367.68
$ws.Range("N136").Value = -2672467.68
